# Weekly crime-stat refresh: shift reporting week forward by one week and
# update the Week-to-Date / 28-Day / Year-to-Date / 2-Year / historical figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header strings: volume/number + reporting week range -------------------
$ws.Range("A8").Value = "Volume 30   Number  14"
$ws.Range("C9").Value = "Report Covering the Week  4/3/2023  Through  4/9/2023"

# --- Plain numeric/percentage cell updates (style unchanged) ----------------
$numericUpdates = @{
    "G15" = 3
    "H15" = -66.666666666666
    "L15" = -33.333333333333
    "M15" = 33.333333333333
    "N15" = 14.285714285714
    "C16" = 11
    "D16" = 5
    "E16" = 120
    "F16" = 25
    "H16" = 8.695652173913
    "I16" = 75
    "J16" = 72
    "K16" = 4.166666666666
    "L16" = 50
    "M16" = 19.047619047619
    "N16" = -64.788732394366
    "C17" = 7
    "E17" = -12.5
    "F17" = 33
    "G17" = 32
    "H17" = 3.125
    "I17" = 101
    "J17" = 104
    "K17" = -2.884615384615
    "L17" = 53.030303030303
    "M17" = -9.009009009009
    "N17" = -0.980392156862
    "C18" = 3
    "D18" = 3
    "E18" = 0
    "G18" = 12
    "H18" = 8.333333333333
    "I18" = 50
    "J18" = 37
    "K18" = 35.135135135135
    "L18" = 56.25
    "M18" = -25.373134328358
    "N18" = -76.851851851851
    "C19" = 9
    "D19" = 11
    "E19" = -18.181818181818
    "F19" = 24
    "G19" = 44
    "H19" = -45.454545454545
    "I19" = 108
    "J19" = 114
    "K19" = -5.263157894736
    "L19" = 52.112676056338
    "M19" = 56.521739130434
    "N19" = 18.681318681318
    "C20" = 5
    "E20" = 66.666666666666
    "F20" = 38
    "G20" = 13
    "H20" = 192.307692307692
    "I20" = 94
    "J20" = 68
    "K20" = 38.235294117647
    "L20" = 203.225806451613
    "M20" = 422.222222222222
    "N20" = 4.444444444444
    "C21" = 35
    "D21" = 30
    "E21" = 16.666666666666
    "F21" = 134
    "G21" = 127
    "H21" = 5.511811023622
    "I21" = 438
    "J21" = 405
    "K21" = 8.148148148148
    "L21" = 66.539923954372
    "M21" = 31.137724550898
    "N21" = -39.752407152682
    "D22" = 1
    "E22" = 0
    "G22" = 3
    "H22" = -66.666666666666
    "I22" = 5
    "J22" = 7
    "K22" = -28.571428571428
    "L22" = 0
    "M22" = 25
    "C24" = 14
    "D24" = 14
    "F24" = 58
    "G24" = 53
    "H24" = 9.433962264150
    "I24" = 196
    "J24" = 234
    "K24" = -16.239316239316
    "L24" = 44.117647058823
    "M24" = 54.330708661417
    "C25" = 11
    "D25" = 11
    "E25" = 0
    "F25" = 38
    "G25" = 58
    "H25" = -34.482758620689
    "I25" = 122
    "J25" = 129
    "K25" = -5.426356589147
    "L25" = 19.607843137254
    "M25" = -10.948905109489
    "D26" = 1
    "E26" = -100
    "J26" = 17
    "K26" = -41.176470588235
    "L26" = -58.333333333333
    "D27" = 3
    "E27" = -66.666666666666
    "F27" = 4
    "G27" = 9
    "H27" = -55.555555555555
    "I27" = 25
    "J27" = 38
    "K27" = -34.210526315789
    "L27" = -16.666666666666
    "M28" = 20
    "N28" = -72.727272727272
    "M29" = 0
    "N29" = -76.190476190476
}
foreach ($ref in $numericUpdates.Keys) {
    $ws.Range($ref).Value = $numericUpdates[$ref]
}

# --- Cells that flip between a suppressed ("0" / "***.*") text marker and ---
# --- an ordinary number: value AND style both change, so after writing the --
# --- value we copy number-format/alignment from a donor cell that already ---
# --- carries the destination style.                                       ---
$styleDonorText   = "N22"   # style 14: right-aligned text marker style
$styleDonorNumber = "D22"   # style 15: right-aligned integer style

$toTextMarker = @{
    "C15" = "0"
    "D15" = "0"
    "E15" = "***.*"
    "D23" = "0"
    "E23" = "***.*"
    "C26" = "0"
    "D28" = "0"
    "E28" = "***.*"
    "D29" = "0"
    "E29" = "***.*"
}
foreach ($ref in $toTextMarker.Keys) {
    $ws.Range($ref).Value = "'" + $toTextMarker[$ref]
    $ws.Range($styleDonorText).Copy()
    $ws.Range($ref).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
}

$toNumber = @{
    "C22" = 1
    "F22" = 1
}
foreach ($ref in $toNumber.Keys) {
    $ws.Range($ref).Value = $toNumber[$ref]
    $ws.Range($styleDonorNumber).Copy()
    $ws.Range($ref).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
}

$excel.CutCopyMode = $false
